$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to Text format for the cells receiving numeric-looking
# strings, so Excel does not silently coerce them into Double values and drop
# meaningful trailing zeros (e.g. "1.140" must stay "1.140", not become 1.14).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.229.16'
$ws.Range("E2").Value = '  +1.07%  '
$ws.Range("D3").Value = '1.868.43'
$ws.Range("E3").Value = '  +3.17%  '
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.45%  '
$ws.Range("D5").Value = '311.89'
$ws.Range("E5").Value = '  +0.86%  '
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("D7").Value = '0.5033'
$ws.Range("E7").Value = '  +0.73%  '
$ws.Range("D8").Value = '0.3904'
$ws.Range("E8").Value = '  +0.37%  '
$ws.Range("D9").Value = '0.09451'
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("D10").Value = '1.140'
$ws.Range("E10").Value = '  +3.68%  '
$ws.Range("D11").Value = '40.84'
$ws.Range("E11").Value = '  +1.04%  '
$ws.Range("D12").Value = '6.444'
$ws.Range("E12").Value = '  +0.49%  '
$ws.Range("D13").Value = '20.94'
$ws.Range("E13").Value = '  +2.18%  '
$ws.Range("D14").Value = '1.868.46'
$ws.Range("E14").Value = '  +3.07%  '
$ws.Range("D15").Value = '1.006'
$ws.Range("E15").Value = '  +0.53%  '
$ws.Range("D16").Value = '7.392'
$ws.Range("E16").Value = '  +1.55%  '
$ws.Range("D17").Value = '0.00001124'
$ws.Range("E17").Value = '  -0.06%  '
$ws.Range("D18").Value = '92.71'
$ws.Range("E18").Value = '  -0.89%  '
$ws.Range("D19").Value = '0.06605'
$ws.Range("E19").Value = '  +0.21%  '
$ws.Range("D20").Value = '17.65'
$ws.Range("E20").Value = '  +2.90%  '
$ws.Range("E21").Value = '  +0.22%  '
$ws.Range("D22").Value = '6.183'
$ws.Range("E22").Value = '  +3.91%  '
$ws.Range("D23").Value = '28.292.13'
$ws.Range("E23").Value = '  +1.13%  '
$ws.Range("D24").Value = '11.27'
$ws.Range("E24").Value = '  +0.82%  '
$ws.Range("D25").Value = '2.293'
$ws.Range("E25").Value = '  +2.02%  '
$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").Value = '3.402'
$ws.Range("E26").Value = '  -0.58%  '
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = '2.569'
$ws.Range("E27").Value = '  +7.03%  '
$ws.Range("B28").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C28").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D28").Value = '2.088.02'
$ws.Range("E28").Value = '  +3.33%  '
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '21.15'
$ws.Range("E29").Value = '  +1.64%  '
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").Value = '158.76'
$ws.Range("E30").Value = '  +0.68%  '
$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D31").Value = '127.31'
$ws.Range("E31").Value = '  -1.15%  '
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").Value = '0.1063'
$ws.Range("E32").Value = '  -1.05%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = '1.063'
$ws.Range("E33").Value = '  +0.75%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = '5.619'
$ws.Range("E34").Value = '  -0.31%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '3.630'
$ws.Range("E35").Value = '  +0.10%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '0.06741'
$ws.Range("E36").Value = '  -1.22%  '
$ws.Range("B37").Value = 'FraxShare'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D37").Value = '9.489'
$ws.Range("E37").Value = '  +5.39%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.02407'
$ws.Range("E38").Value = '  +3.84%  '
$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").Value = '0.2188'
$ws.Range("E39").Value = '  +1.80%  '
$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").Value = '11.49'
$ws.Range("E40").Value = '  +1.02%  '
$ws.Range("D41").Value = '4.998'
$ws.Range("E41").Value = '  +0.95%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '0.6346'
$ws.Range("E42").Value = '  +1.23%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '1.184'
$ws.Range("E43").Value = '  +3.38%  '
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").Value = '1.002'
$ws.Range("E44").Value = '  +0.25%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '13.56'
$ws.Range("E45").Value = '  +3.72%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '0.5980'
$ws.Range("E46").Value = '  +1.07%  '
$ws.Range("B47").Value = 'WEMIXTOKEN'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").Value = '1.278'
$ws.Range("E47").Value = '  -1.09%  '
$ws.Range("B48").Value = 'PancakeSwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D48").Value = '3.659'
$ws.Range("E48").Value = '  -0.88%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '2.000'
$ws.Range("E49").Value = '  +2.01%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = '123.23'
$ws.Range("E50").Value = '  -0.79%  '
$ws.Range("B51").Value = 'EOS'
$ws.Range("C51").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D51").Value = '1.197'
$ws.Range("E51").Value = '  +1.44%  '
